$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A27").Value = "18-10-2025"
$ws.Range("B27").Value = "The price of gold in India today is ₹13,086 per gram for 24 karat gold, ₹11,995 per gram for 22 karat gold and ₹9,814 per gram for 18 karat gold (also called 999 gold)."

$ws.Range("A27").Borders.LineStyle = 1
$ws.Range("B27").Borders.LineStyle = 1
$ws.Range("B27").WrapText = $true
